$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.89
$ws.Range("C2").Value = 0.61
$ws.Range("F2").Value = -0.13
$ws.Range("G2").Value = 1.12
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0.01

# Row 3
$ws.Range("B3").Value = 0.32
$ws.Range("C3").Value = 1.84
$ws.Range("F3").Value = 1.06
$ws.Range("G3").Value = 1.61
$ws.Range("H3").Value = 0.9
$ws.Range("I3").Value = 0.37

# Row 4
$ws.Range("B4").Value = 1.06
$ws.Range("C4").Value = 1.61
$ws.Range("F4").Value = 1.06
$ws.Range("G4").Value = 1.61
$ws.Range("I4").Value = 0.41
